# Insert a new data row at row 205 (pushes existing rows 205..278 down to 206..279)
# and populate it with the new record described by the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("205:205").Insert()

$ws.Range("A205").Value = 4
$ws.Range("B205").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C205").Value = "Los Lagos"
$ws.Range("D205").Value = 44876
$ws.Range("E205").Value = 10
$ws.Range("F205").Value = 100112039
$ws.Range("G205").Value = "Ciboulette"
$ws.Range("H205").Value = "Sin especificar"
$ws.Range("I205").Value = "Primera"
$ws.Range("J205").Value = 240
$ws.Range("K205").Value = 2500
$ws.Range("L205").Value = 2500
$ws.Range("M205").Value = 2500
$ws.Range("N205").Value = "$/docena de atados"
$ws.Range("O205").Value = "Región Metropolitana"
$ws.Range("P205").Value = 833
$ws.Range("Q205").Value = 3
$ws.Range("R205").Value = "Hortaliza"
